$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (stored as serial number 45190)
# that was bumped by 2 days to serial number 45192 for every data row
# (rows 2 through 205).
$newDate = Get-Date -Year 2023 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
$ws.Range("C2:C$lastRow").Value = $newDate
